# Applies the "feat(form): ui additions and key changes" edit to the
# Final Survey Report document.
#
# Word constants used below (COM-interop convention):
#   wdReplaceAll = 2
#   wdFindContinue = 1 (wrap mode used positionally, not referenced by name)

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- Header block: team ref + date ---------------------------------------
Replace-Text "DEL/NIA/10229/FSR" "123456789"
Replace-Text "29th July 2025" "06/05/2025"

# --- Addressee block --------------------------------------------------
Replace-Text "The New India Assurance Co. Ltd.," "Nike,"
Replace-Text "Centralized Claims Hub MRO-III, New India Centre," "D-50a, 1st floor, Pandav Nagar,"
Replace-Text "3rd Floor, 17-A, Cooperage Road, Mumbai – 400 001." "Delhi-110092"
Replace-Text "Kind Attn.: Ms. Swati Tilak/ Ms. Indra Ayer." "Kind Attn.: Sonal Singh"

# --- Policy numbers (one zero removed from each long number) -------------
Replace-Text "Cellular Network Policy No. 121200462426000000001" "Cellular Network Policy No. 12120046242600000001"
Replace-Text "Burglary & Theft. 121200462426000000002" "Burglary & Theft. 12120046242600000002"

# --- Claim no. block: spaces -> colon -------------------------------------
Replace-Text "Insurer Claim No.    Details Awaited" "Insurer Claim No.: Details Awaited"
Replace-Text "Edme Control No.    MIS " "Edme Control No.: MIS "

# --- Insured address: merge the underlined "Budhkar" run into plain text,
#     and drop the trailing period ----------------------------------------
Replace-Text "Pandurang Budhkar Marg, Worli, Mumbai – 400 030." "Pandurang Budhkar Marg, Worli, Mumbai – 400 030"

# --- Address of loss: expand the truncated text ---------------------------
Replace-Text "Shri. Ram Kishan S/o Shri. Hukami. Khewat/Khata No.- 261/330. Rect" `
    "Shri. Ram Kishan S/o Shri. Hukami, Khewat/Khata No.- 261/330, Rect No.- 63, Killa No.- 4/1/3, Village & Post Office-Bagpur, Tehsil & District-Palwal, Haryana, India. Pin Code-121102"

# --- Person contacted table: name + phone ---------------------------------
Replace-Text "Nitansh" "Gomez"
Replace-Text "9813017817" "1100110011"

# --- Person contacted table: widen the "Mr./Mrs." column and narrow the
#     name/number value column (45% -> 25%) for both rows -----------------
$tbl = $d.Tables.Item(3)
Write-Output ("table count: " + $d.Tables.Count)
